# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new labels, styled like the rest of the header row ---
$ws.Range("A1").Copy() | Out-Null
$ws.Range("AD1:AF1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (2-40): same season record (101-59-1) for every player ---
$ws.Range("AD2:AD40").Value = 101
$ws.Range("AE2:AE40").Value = 59
$ws.Range("AF2:AF40").Value = 1
